# Reduce mannings n and increase slope. Go faster to stop the piling effect
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("units_m_day")
$ws.Activate()

# mannings_n (column V) for rows 2-11 reduced from 0.1 to 0.08
$ws.Range("V2:V11").Value = 0.08

# Match the author's final selection on the sheet
$ws.Range("T12").Select()
